# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price-looking D values are forced to Text ("@") before the write so Excel
# doesn't auto-coerce them to numbers, then the style is reset to "Normal"
# so no stray number-format style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '41.771.53'
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").Value = '2.267.35'
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("E7").Value = '  +1.50%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("E12").Value = '  +0.22%  '

$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").Value = '2.614.38'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.87%  '

$ws.Range("D17").Value = '2.261.11'
$ws.Range("E17").Value = '  +1.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.766'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '

$ws.Range("D19").Value = '41.681.18'
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.51%  '

$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.01'
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  +1.18%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  +2.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.43%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  -5.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.38%  '

$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0742'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.43%  '

$ws.Range("E36").Value = '  -1.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.31%  '

$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("E39").Value = '  +1.26%  '

$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").Value = '2.025.15'
$ws.Range("E43").Value = '  -3.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0279'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '72.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.39%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.81%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
